$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.082.91"
$ws.Range("E2").Value = "  -1.71%  "
$ws.Range("D3").Value = "3.736.86"
$ws.Range("E3").Value = "  -1.23%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "621.35"
$ws.Range("E5").Value = "  +0.19%  "
$ws.Range("D6").Value = "180.79"
$ws.Range("E6").Value = "  -0.71%  "
$ws.Range("D7").Value = "3.734.84"
$ws.Range("E7").Value = "  -1.37%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("E10").Value = "  +1.55%  "
$ws.Range("E11").Value = "  -4.80%  "
$ws.Range("D12").Value = "0.486"
$ws.Range("E12").Value = "  -3.56%  "
$ws.Range("D13").Value = "40.75"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "0.0000259"
$ws.Range("E14").Value = "  +0.71%  "
$ws.Range("D15").Value = "4.357.11"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "3.734.89"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "70.094.07"
$ws.Range("E17").Value = "  -1.75%  "
$ws.Range("E18").Value = "  -1.75%  "
$ws.Range("D19").Value = "7.61"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("D20").Value = "16.77"
$ws.Range("E20").Value = "  -1.26%  "
$ws.Range("D21").Value = "506.11"
$ws.Range("E21").Value = "  -2.92%  "
$ws.Range("E22").Value = "  -1.16%  "
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").Value = "2.55"
$ws.Range("E24").Value = "  -0.08%  "
$ws.Range("D25").Value = "86.67"
$ws.Range("E25").Value = "  -2.68%  "
$ws.Range("D26").Value = "'11.50"
$ws.Range("E26").Value = "  +2.57%  "
$ws.Range("D27").Value = "13.12"
$ws.Range("E27").Value = "  -3.76%  "
$ws.Range("E28").Value = "  +21.09%  "
$ws.Range("E29").Value = "  -0.14%  "
$ws.Range("E30").Value = "  -2.59%  "
$ws.Range("E31").Value = "  +0.69%  "
$ws.Range("D32").Value = "7.94"
$ws.Range("E32").Value = "  -2.73%  "
$ws.Range("D33").Value = "31.22"
$ws.Range("E33").Value = "  -3.73%  "
$ws.Range("E34").Value = "  -1.45%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  -6.64%  "
$ws.Range("D41").Value = "50.22"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("D42").Value = "45.54"
$ws.Range("E42").Value = "  +1.20%  "
$ws.Range("D43").Value = "434.95"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("D44").Value = "2.88"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  -2.15%  "
$ws.Range("D46").Value = "3.004.34"
$ws.Range("E46").Value = "  -5.50%  "
$ws.Range("E47").Value = "  -0.93%  "
$ws.Range("D48").Value = "27.56"
$ws.Range("E48").Value = "  -2.36%  "
$ws.Range("D50").Value = "137.28"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  +0.45%  "
